$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AMOUNT column (B2:B6) to 5000 with a 2-decimal number format
$ws.Range("B2:B6").Value = 5000
$ws.Range("B2:B6").NumberFormat = "0.00"

# Update PAYMENT_METHOD column (D2:D6) to CREDIT_CARD
$ws.Range("D2:D6").Value = "CREDIT_CARD"

# Update STATUS column (E2:E6) to REFUNDED
$ws.Range("E2:E6").Value = "REFUNDED"

# Update selection
$ws.Range("G5").Select()
